$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header values for the additional year columns (P4, Q4, R4) and data (P5, Q5, R5)
$ws.Range("P4").Value = 2019
$ws.Range("Q4").Value = 2020
$ws.Range("R4").Value = 2021

$ws.Range("P5").Value = 12.9
$ws.Range("Q5").Value = 15.2
$ws.Range("R5").Value = 10.4

# Copy style from the existing thick-bottom-border row (row 3) to the newly used Q3/R3 cells
$ws.Range("O3").Copy() | Out-Null
$ws.Range("P3:R3").PasteSpecial(-4122) | Out-Null

# Copy style (font+border+alignment) from O4 (years row) to P4:R4
$ws.Range("O4").Copy() | Out-Null
$ws.Range("P4:R4").PasteSpecial(-4122) | Out-Null

# Copy style from O5 (data row) to P5:R5
$ws.Range("O5").Copy() | Out-Null
$ws.Range("P5:R5").PasteSpecial(-4122) | Out-Null

$ws.Range("S3").Select() | Out-Null
